$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B8").Value = 45692
$ws.Range("B9").Value = 45694

$ws.Range("H9").Value = "HW 1"
$ws.Range("H8").ClearContents()

$ws.Range("F16").Select() | Out-Null

Write-Host "done"
